# Updates the "Estado de Cuenta" workbook:
#  - Removes the first worker's (LUIS MIGUEL ARRIETA LOPEZ / CC 9094130) 5 rows
#    of period data (rows 16-20), shifting the remaining worker's rows and the
#    signature footer up.
#  - Reverses the remaining worker's (LUISA MIGELDY ARRIETA BUSTOS) period rows
#    so periods run in ascending order (1804 .. 1906) instead of descending.
#  - Refreshes the summary header cells (total mora, worker count, period count).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the 5 rows belonging to the worker who no longer appears in the
# statement. This also shifts the trailing signature rows (40/41 -> 35/36)
# and updates merged cells / dimension automatically.
$ws.Range("B16:J20").EntireRow.Delete()

# The remaining worker's 15 period rows (now at 16..30) were in descending
# period order (1906 .. 1804); put them back in ascending order (1804 .. 1906).
$periodCount = 15
$firstRow = 16
$periods = @()
for ($i = 0; $i -lt $periodCount; $i++) {
    $periods += $ws.Range("E" + ($firstRow + $i)).Value()
}
for ($i = 0; $i -lt $periodCount; $i++) {
    $ws.Range("E" + ($firstRow + $i)).Value = $periods[$periodCount - 1 - $i]
}

# Refresh the summary figures at the top of the statement.
$ws.Range("E11").Value = 468735   # VALOR MORA total
$ws.Range("C13").Value = 1        # Cant. Trabajadores
$ws.Range("F13").Value = 15       # Cant. Periodos
